# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted as row 32 of the data table
# (Feria Lagunitas de Puerto Montt - Apio), pushing every existing row
# from 32 downward to the next row. The new row carries:
#   Fecha (D) = 2021-08-06 (serial 44414), Volumen (J) = 40
# with the remaining columns matching the values already used by the
# (pre-insert) row 32 entry (Primera / 12000 / 12000 / 12000 / $/docena
# de matas / Región de Coquimbo / 2000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; existing rows 32..108 shift down to 33..109.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly observation.
$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("D32").Value = 44414
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 100112017
$ws.Range("G32").Value = "Apio"
$ws.Range("H32").Value = "Americana (o)"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 40
$ws.Range("K32").Value = 12000
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = 12000
$ws.Range("N32").Value = "`$/docena de matas"
$ws.Range("O32").Value = "Región de Coquimbo"
$ws.Range("P32").Value = 2000
$ws.Range("Q32").Value = 6
$ws.Range("R32").Value = "Hortaliza"
